$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = 'presidenital'
$ws.Cells.Item(2,3).Value = 'presidential'
$ws.Cells.Item(2,4).Value = 'He oversaw the transformation of Turkey''s parliamentary system into a presidenital system, introducing term limits and expanding executive powers, and Turkey''s migrant crisis.'
$ws.Cells.Item(2,5).Value = 'He oversaw the transformation of Turkey''s parliamentary system into a presidential system, introducing term limits and expanding executive powers, and Turkey''s migrant crisis.'
$ws.Cells.Item(3,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(3,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(3,4).Value = 'Many psychologists are involved in some kind of therapeutic role, practicing psychotherapy in clinical, counseling, or school settings.'
$ws.Cells.Item(3,5).Value = 'Many psychologists are involved in some kind of therapeutic role, practicing psychotherapy in clinical, counseling, or school settings.'
$ws.Cells.Item(4,4).Value = 'This was the ninth opened McDonald''s restaurant overall, although this location was demolished in 1984 in response to the San Ysidro McDonald''s massacre.'
$ws.Cells.Item(4,5).Value = 'This was the ninth opened McDonald''s restaurant overall, although this location was demolished in 1984 in response to the San Ysidro McDonald''s massacre.'
$ws.Cells.Item(5,2).Value = 'preiods'
$ws.Cells.Item(5,3).Value = 'periods'
$ws.Cells.Item(5,4).Value = 'Thereafter, preiods of civil war and Seljuk incursion resulted in the loss of most of Asia Minor.'
$ws.Cells.Item(5,5).Value = 'Thereafter, periods of civil war and Seljuk incursion resulted in the loss of most of Asia Minor.'
$ws.Cells.Item(6,2).Value = 'Ofxord'
$ws.Cells.Item(6,3).Value = 'Oxford'
$ws.Cells.Item(6,4).Value = 'Ofxord does not have a main campus.'
$ws.Cells.Item(6,5).Value = 'Oxford does not have a main campus.'
$ws.Cells.Item(7,2).Value = 'sapre'
$ws.Cells.Item(7,3).Value = 'spare'
$ws.Cells.Item(7,4).Value = 'Almost every kind of product can be found in the international market, for example: food, clothes, sapre parts, oil, jewellery, wine, stocks, currencies, and water.'
$ws.Cells.Item(7,5).Value = 'Almost every kind of product can be found in the international market, for example: food, clothes, spare parts, oil, jewellery, wine, stocks, currencies, and water.'
$ws.Cells.Item(8,2).Value = 'appiles'
$ws.Cells.Item(8,3).Value = 'applies'
$ws.Cells.Item(8,4).Value = 'The judiciary (also known as the judicial system, judicature, judicial branch, judiciative branch, and court or judiciary system) is the system of courts that adjudicates legal disputes/disagreements and interprets, defends, and appiles the law in legal cases.'
$ws.Cells.Item(8,5).Value = 'The judiciary (also known as the judicial system, judicature, judicial branch, judiciative branch, and court or judiciary system) is the system of courts that adjudicates legal disputes/disagreements and interprets, defends, and applies the law in legal cases.'
$ws.Cells.Item(9,2).Value = 'dritfed'
$ws.Cells.Item(9,3).Value = 'drifted'
$ws.Cells.Item(9,4).Value = 'Later he dritfed into ill-health and solitude.'
$ws.Cells.Item(9,5).Value = 'Later he drifted into ill-health and solitude.'
$ws.Cells.Item(10,2).Value = 'intnetions'
$ws.Cells.Item(10,3).Value = 'intentions'
$ws.Cells.Item(10,4).Value = 'The Truman Doctrine in 1947 enunciated American intnetions to guarantee the security of Turkey and Greece during the Cold War, and resulted in large-scale U.S.'
$ws.Cells.Item(10,5).Value = 'The Truman Doctrine in 1947 enunciated American intentions to guarantee the security of Turkey and Greece during the Cold War, and resulted in large-scale U.S.'
$ws.Cells.Item(11,2).Value = 'reevnue'
$ws.Cells.Item(11,3).Value = 'revenue'
$ws.Cells.Item(11,4).Value = 'In 2023, YouTube''s advertising reevnue totaled $31.7 billion, a 2% increase from the $31.1 billion reported in 2022.'
$ws.Cells.Item(11,5).Value = 'In 2023, YouTube''s advertising revenue totaled $31.7 billion, a 2% increase from the $31.1 billion reported in 2022.'
$ws.Cells.Item(12,2).Value = 'internaitonally'
$ws.Cells.Item(12,3).Value = 'internationally'
$ws.Cells.Item(12,4).Value = 'It is available internaitonally in multiple languages.'
$ws.Cells.Item(12,5).Value = 'It is available internationally in multiple languages.'
$ws.Cells.Item(13,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(13,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(13,4).Value = 'Historically, parliaments included various kinds of deliberative, consultative, and judicial assemblies.'
$ws.Cells.Item(13,5).Value = 'Historically, parliaments included various kinds of deliberative, consultative, and judicial assemblies.'
$ws.Cells.Item(14,2).Value = 'Presidentail'
$ws.Cells.Item(14,3).Value = 'Presidential'
$ws.Cells.Item(14,4).Value = 'In 2018, Presley was posthumously awarded the Presidentail Medal of Freedom.'
$ws.Cells.Item(14,5).Value = 'In 2018, Presley was posthumously awarded the Presidential Medal of Freedom.'
$ws.Cells.Item(15,2).Value = 'Truks'
$ws.Cells.Item(15,3).Value = 'Turks'
$ws.Cells.Item(15,4).Value = 'The Seljuk Truks began migrating into Anatolia in the 11th century, starting the Turkification process.'
$ws.Cells.Item(15,5).Value = 'The Seljuk Turks began migrating into Anatolia in the 11th century, starting the Turkification process.'
$ws.Cells.Item(16,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(16,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(16,4).Value = 'He began his music career in 1954 at Sun Records with producer Sam Phillips, who wanted to bring the sound of African-American music to a wider audience.'
$ws.Cells.Item(16,5).Value = 'He began his music career in 1954 at Sun Records with producer Sam Phillips, who wanted to bring the sound of African-American music to a wider audience.'
$ws.Cells.Item(17,2).Value = 'wihch'
$ws.Cells.Item(17,3).Value = 'which'
$ws.Cells.Item(17,4).Value = 'It is predated in use by the Ancient Greek word τέχνη (tékhnē), used to mean ''knowledge of how to make things'', wihch encompassed activities like architecture.'
$ws.Cells.Item(17,5).Value = 'It is predated in use by the Ancient Greek word τέχνη (tékhnē), used to mean ''knowledge of how to make things'', which encompassed activities like architecture.'
$ws.Cells.Item(18,2).Value = 'riisng'
$ws.Cells.Item(18,3).Value = 'rising'
$ws.Cells.Item(18,4).Value = 'With riisng nationalism, a number of new states emerged in the Balkans.'
$ws.Cells.Item(18,5).Value = 'With rising nationalism, a number of new states emerged in the Balkans.'
$ws.Cells.Item(19,2).Value = 'ofifce'
$ws.Cells.Item(19,3).Value = 'office'
$ws.Cells.Item(19,4).Value = 'The ofifce was later renamed to Minister-President of the Austrian Empire and remained from there on until the dissolution of Austria-Hungary.'
$ws.Cells.Item(19,5).Value = 'The office was later renamed to Minister-President of the Austrian Empire and remained from there on until the dissolution of Austria-Hungary.'
$ws.Cells.Item(20,2).Value = 'conudct'
$ws.Cells.Item(20,3).Value = 'conduct'
$ws.Cells.Item(20,4).Value = 'It proscribes conudct perceived as threatening, harmful, or otherwise endangering to the property, health, safety, and welfare of people inclusive of one''s self.'
$ws.Cells.Item(20,5).Value = 'It proscribes conduct perceived as threatening, harmful, or otherwise endangering to the property, health, safety, and welfare of people inclusive of one''s self.'
$ws.Cells.Item(21,2).Value = 'glboal'
$ws.Cells.Item(21,3).Value = 'global'
$ws.Cells.Item(21,4).Value = 'Poorer communities are responsible for a small share of glboal emissions, yet have the least ability to adapt and are most vulnerable to climate change.'
$ws.Cells.Item(21,5).Value = 'Poorer communities are responsible for a small share of global emissions, yet have the least ability to adapt and are most vulnerable to climate change.'
$ws.Cells.Item(22,2).Value = 'acnient'
$ws.Cells.Item(22,3).Value = 'ancient'
$ws.Cells.Item(22,4).Value = 'Classical Greek culture, especially philosophy, had a powerful influence on acnient Rome, which carried a version of it throughout the Mediterranean and much of Europe.'
$ws.Cells.Item(22,5).Value = 'Classical Greek culture, especially philosophy, had a powerful influence on ancient Rome, which carried a version of it throughout the Mediterranean and much of Europe.'
$ws.Cells.Item(23,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(23,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(23,4).Value = 'Alternatively, this document can simply be referred to as a degree certificate or graduation certificate, or as a parchment.'
$ws.Cells.Item(23,5).Value = 'Alternatively, this document can simply be referred to as a degree certificate or graduation certificate, or as a parchment.'
$ws.Cells.Item(24,2).Value = 'ytpes'
$ws.Cells.Item(24,3).Value = 'types'
$ws.Cells.Item(24,4).Value = 'The most widely used renewable energy ytpes are solar energy, wind power, and hydropower.'
$ws.Cells.Item(24,5).Value = 'The most widely used renewable energy types are solar energy, wind power, and hydropower.'
$ws.Cells.Item(25,2).Value = 'coheernt'
$ws.Cells.Item(25,3).Value = 'coherent'
$ws.Cells.Item(25,4).Value = 'Historians strive to integrate the perspectives of several sources to develop a coheernt narrative.'
$ws.Cells.Item(25,5).Value = 'Historians strive to integrate the perspectives of several sources to develop a coherent narrative.'
$ws.Cells.Item(26,2).Value = 'rpogram'
$ws.Cells.Item(26,3).Value = 'program'
$ws.Cells.Item(26,4).Value = 'Runs a program with arguments for each candidate file, selecting files for which the rpogram results in exit status 0.'
$ws.Cells.Item(26,5).Value = 'Runs a program with arguments for each candidate file, selecting files for which the program results in exit status 0.'
$ws.Cells.Item(27,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(27,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(27,4).Value = 'The Italian Renaissance concluded in 1527 when Holy Roman Emperor Charles V launched an assault on Rome during the war of the League of Cognac.'
$ws.Cells.Item(27,5).Value = 'The Italian Renaissance concluded in 1527 when Holy Roman Emperor Charles V launched an assault on Rome during the war of the League of Cognac.'
$ws.Cells.Item(28,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(28,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(28,4).Value = 'One of the most revered figures in the history of Western music, his works rank among the most performed of the classical music repertoire and span the transition from the Classical to the Romantic era.'
$ws.Cells.Item(28,5).Value = 'One of the most revered figures in the history of Western music, his works rank among the most performed of the classical music repertoire and span the transition from the Classical to the Romantic era.'
$ws.Cells.Item(29,2).Value = 'acpital'
$ws.Cells.Item(29,3).Value = 'capital'
$ws.Cells.Item(29,4).Value = 'Kenya''s acpital and largest city is Nairobi.'
$ws.Cells.Item(29,5).Value = 'Kenya''s capital and largest city is Nairobi.'
$ws.Cells.Item(30,2).Value = 'populatino'
$ws.Cells.Item(30,3).Value = 'population'
$ws.Cells.Item(30,4).Value = 'The Late Middle Ages was marked by difficulties and calamities, including famine, plague, and war, which significantly diminished the populatino of Europe; between 1347 and 1350, the Black Death killed about a third of Europeans.'
$ws.Cells.Item(30,5).Value = 'The Late Middle Ages was marked by difficulties and calamities, including famine, plague, and war, which significantly diminished the population of Europe; between 1347 and 1350, the Black Death killed about a third of Europeans.'
$ws.Cells.Item(31,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(31,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(31,4).Value = 'It starts searching at a directory in a hierarchical structure and recursively traverses the tree although can be limited to a maximum number of levels.'
$ws.Cells.Item(31,5).Value = 'It starts searching at a directory in a hierarchical structure and recursively traverses the tree although can be limited to a maximum number of levels.'
$ws.Cells.Item(32,2).Value = 'activiites'
$ws.Cells.Item(32,3).Value = 'activities'
$ws.Cells.Item(32,4).Value = 'These activiites include painting, sculpting, music, theatre, literature, and more.'
$ws.Cells.Item(32,5).Value = 'These activities include painting, sculpting, music, theatre, literature, and more.'
$ws.Cells.Item(33,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(33,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(33,4).Value = 'Its remaining territories were progressively annexed by the Ottomans in a series of wars fought in the 14th and 15th centuries.'
$ws.Cells.Item(33,5).Value = 'Its remaining territories were progressively annexed by the Ottomans in a series of wars fought in the 14th and 15th centuries.'
$ws.Cells.Item(34,2).Value = 'chattle'
$ws.Cells.Item(34,3).Value = 'chattel'
$ws.Cells.Item(34,4).Value = 'Constitution, which in 1865 abolished chattle slavery.'
$ws.Cells.Item(34,5).Value = 'Constitution, which in 1865 abolished chattel slavery.'
$ws.Cells.Item(35,2).Value = 'pAplied'
$ws.Cells.Item(35,3).Value = 'Applied'
$ws.Cells.Item(35,4).Value = 'pAplied sociological research may be applied directly to social policy and welfare, whereas theoretical approaches may focus on the understanding of social processes and phenomenological method.'
$ws.Cells.Item(35,5).Value = 'Applied sociological research may be applied directly to social policy and welfare, whereas theoretical approaches may focus on the understanding of social processes and phenomenological method.'
$ws.Cells.Item(36,2).Value = 'aHrry'
$ws.Cells.Item(36,3).Value = 'Harry'
$ws.Cells.Item(36,4).Value = 'In 1948, President aHrry S.'
$ws.Cells.Item(36,5).Value = 'In 1948, President Harry S.'
$ws.Cells.Item(37,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(37,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(37,4).Value = 'The modern word alchemy in turn is derived from the Arabic word al-kīmīā (الكیمیاء).'
$ws.Cells.Item(37,5).Value = 'The modern word alchemy in turn is derived from the Arabic word al-kīmīā (الكیمیاء).'
$ws.Cells.Item(38,4).Value = 'CPJ staff applies strict criteria for each case; researchers independently investigate and verify the circumstances behind each death or imprisonment.'
$ws.Cells.Item(38,5).Value = 'CPJ staff applies strict criteria for each case; researchers independently investigate and verify the circumstances behind each death or imprisonment.'
$ws.Cells.Item(39,2).Value = 'farnchised'
$ws.Cells.Item(39,3).Value = 'franchised'
$ws.Cells.Item(39,4).Value = 'They soon farnchised the company.'
$ws.Cells.Item(39,5).Value = 'They soon franchised the company.'
$ws.Cells.Item(40,2).Value = 'obsrevable'
$ws.Cells.Item(40,3).Value = 'observable'
$ws.Cells.Item(40,4).Value = 'Nevertheless, a scholarly tradition obsrevable in the 16th century claimed a far earlier 7th-century foundation by Archbishop Theodore of Tarsus (668–690).'
$ws.Cells.Item(40,5).Value = 'Nevertheless, a scholarly tradition observable in the 16th century claimed a far earlier 7th-century foundation by Archbishop Theodore of Tarsus (668–690).'
$ws.Cells.Item(41,2).Value = 'ubran'
$ws.Cells.Item(41,3).Value = 'urban'
$ws.Cells.Item(41,4).Value = 'The ubran population of the city is 330,836 (Ortahisar), with a metropolitan population of 822,270.'
$ws.Cells.Item(41,5).Value = 'The urban population of the city is 330,836 (Ortahisar), with a metropolitan population of 822,270.'
$ws.Cells.Item(42,2).Value = 'calssical'
$ws.Cells.Item(42,3).Value = 'classical'
$ws.Cells.Item(42,4).Value = 'The terminology changed in the 15th century as the renewed interest in the writings of Ancient Rome caused writers to prefer calssical terminology.'
$ws.Cells.Item(42,5).Value = 'The terminology changed in the 15th century as the renewed interest in the writings of Ancient Rome caused writers to prefer classical terminology.'
$ws.Cells.Item(43,2).Value = 'aobut'
$ws.Cells.Item(43,3).Value = 'about'
$ws.Cells.Item(43,4).Value = 'Ethical concerns have been raised aobut AI''s long-term effects and potential existential risks, prompting discussions about regulatory policies to ensure the safety and benefits of the technology.'
$ws.Cells.Item(43,5).Value = 'Ethical concerns have been raised about AI''s long-term effects and potential existential risks, prompting discussions about regulatory policies to ensure the safety and benefits of the technology.'
$ws.Cells.Item(44,2).Value = 'Tukrish'
$ws.Cells.Item(44,3).Value = 'Turkish'
$ws.Cells.Item(44,4).Value = 'The city is Turkey''s biggest international sea resort on the Tukrish Riviera.'
$ws.Cells.Item(44,5).Value = 'The city is Turkey''s biggest international sea resort on the Turkish Riviera.'
$ws.Cells.Item(45,2).Value = 'restaurnat'
$ws.Cells.Item(45,3).Value = 'restaurant'
$ws.Cells.Item(45,4).Value = 'Today, McDonald''s has more than 40,000 restaurnat locations worldwide, with around one-quarter in the US.'
$ws.Cells.Item(45,5).Value = 'Today, McDonald''s has more than 40,000 restaurant locations worldwide, with around one-quarter in the US.'
$ws.Cells.Item(46,2).Value = 'economci'
$ws.Cells.Item(46,3).Value = 'economic'
$ws.Cells.Item(46,4).Value = 'Climate change threatens people with increased flooding, extreme heat, increased food and water scarcity, more disease, and economci loss.'
$ws.Cells.Item(46,5).Value = 'Climate change threatens people with increased flooding, extreme heat, increased food and water scarcity, more disease, and economic loss.'
$ws.Cells.Item(47,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(47,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(47,4).Value = 'On average, the Supreme Court receives about 7,000 petitions for writs of certiorari each year, but only grants about 80.'
$ws.Cells.Item(47,5).Value = 'On average, the Supreme Court receives about 7,000 petitions for writs of certiorari each year, but only grants about 80.'
$ws.Cells.Item(48,2).Value = 'judicail'
$ws.Cells.Item(48,3).Value = 'judicial'
$ws.Cells.Item(48,4).Value = 'The first phase, In Iure, was the judicail process.'
$ws.Cells.Item(48,5).Value = 'The first phase, In Iure, was the judicial process.'
$ws.Cells.Item(49,2).Value = 'bureaucrcay'
$ws.Cells.Item(49,3).Value = 'bureaucracy'
$ws.Cells.Item(49,4).Value = 'Communism''s decline has been attributed to economic inefficiency and to authoritarianism and bureaucrcay within Communist governments.'
$ws.Cells.Item(49,5).Value = 'Communism''s decline has been attributed to economic inefficiency and to authoritarianism and bureaucracy within Communist governments.'
$ws.Cells.Item(50,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(50,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(50,4).Value = 'In 2000, the International Monetary Fund (IMF) identified four basic aspects of globalization: trade and transactions, capital and investment movements, migration and movement of people, and the dissemination of knowledge.'
$ws.Cells.Item(50,5).Value = 'In 2000, the International Monetary Fund (IMF) identified four basic aspects of globalization: trade and transactions, capital and investment movements, migration and movement of people, and the dissemination of knowledge.'
$ws.Cells.Item(51,2).Value = 'Amreican'
$ws.Cells.Item(51,3).Value = 'American'
$ws.Cells.Item(51,4).Value = 'In addition, 574 Native American tribes have sovereignty rights, and there are 326 Native Amreican reservations.'
$ws.Cells.Item(51,5).Value = 'In addition, 574 Native American tribes have sovereignty rights, and there are 326 Native American reservations.'
$ws.Cells.Item(52,2).Value = 'ccyling'
$ws.Cells.Item(52,3).Value = 'cycling'
$ws.Cells.Item(52,4).Value = 'Medicina has a football and a basketball team, both playing in lower leagues, as well as a ccyling team and a cycling development center for youth.'
$ws.Cells.Item(52,5).Value = 'Medicina has a football and a basketball team, both playing in lower leagues, as well as a cycling team and a cycling development center for youth.'
$ws.Cells.Item(53,2).Value = 'stereotyipcally'
$ws.Cells.Item(53,3).Value = 'stereotypically'
$ws.Cells.Item(53,4).Value = 'The chemistry laboratory stereotyipcally uses various forms of laboratory glassware.'
$ws.Cells.Item(53,5).Value = 'The chemistry laboratory stereotypically uses various forms of laboratory glassware.'
$ws.Cells.Item(54,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(54,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(54,4).Value = 'The historic centre of Istanbul is a UNESCO World Heritage Site.'
$ws.Cells.Item(54,5).Value = 'The historic centre of Istanbul is a UNESCO World Heritage Site.'
$ws.Cells.Item(55,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(55,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(55,4).Value = 'While the position does not have a formal term limit, incumbent secretary-generals have avoided seeking a third term since the 1981 selection, when China cast a record 16 vetoes against a third term for Kurt Waldheim.'
$ws.Cells.Item(55,5).Value = 'While the position does not have a formal term limit, incumbent secretary-generals have avoided seeking a third term since the 1981 selection, when China cast a record 16 vetoes against a third term for Kurt Waldheim.'
$ws.Cells.Item(56,2).Value = 'conrtol'
$ws.Cells.Item(56,3).Value = 'control'
$ws.Cells.Item(56,4).Value = 'Turkey''s concept of laiklik ("laicism") calls for the separation of state and religion, but also describes the state''s stance as one of "active neutrality", which involves state conrtol and legal regulation of religion.'
$ws.Cells.Item(56,5).Value = 'Turkey''s concept of laiklik ("laicism") calls for the separation of state and religion, but also describes the state''s stance as one of "active neutrality", which involves state control and legal regulation of religion.'
$ws.Cells.Item(57,2).Value = 'insttiutions'
$ws.Cells.Item(57,3).Value = 'institutions'
$ws.Cells.Item(57,4).Value = 'All of the colleges are self-governing insttiutions within the university, managing their own personnel and policies, and all students are required to have a college affiliation within the university.'
$ws.Cells.Item(57,5).Value = 'All of the colleges are self-governing institutions within the university, managing their own personnel and policies, and all students are required to have a college affiliation within the university.'
$ws.Cells.Item(58,2).Value = 'estabilshed'
$ws.Cells.Item(58,3).Value = 'established'
$ws.Cells.Item(58,4).Value = 'Darwin''s work estabilshed evolutionary descent with modification as the dominant scientific explanation of natural diversification.'
$ws.Cells.Item(58,5).Value = 'Darwin''s work established evolutionary descent with modification as the dominant scientific explanation of natural diversification.'
$ws.Cells.Item(59,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(59,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(59,4).Value = 'The January 1975 issue of Popular Electronics featured Micro Instrumentation and Telemetry Systems''s (MITS) Altair 8800 microcomputer, which inspired Allen to suggest that they could program a BASIC interpreter for the device.'
$ws.Cells.Item(59,5).Value = 'The January 1975 issue of Popular Electronics featured Micro Instrumentation and Telemetry Systems''s (MITS) Altair 8800 microcomputer, which inspired Allen to suggest that they could program a BASIC interpreter for the device.'
$ws.Cells.Item(60,2).Value = 'cerdited'
$ws.Cells.Item(60,3).Value = 'credited'
$ws.Cells.Item(60,4).Value = 'Often regarded as one of the most influential entertainment brands in history, Disney is cerdited with revolutionizing the animation industry.'
$ws.Cells.Item(60,5).Value = 'Often regarded as one of the most influential entertainment brands in history, Disney is credited with revolutionizing the animation industry.'
$ws.Cells.Item(61,2).Value = 'seaprates'
$ws.Cells.Item(61,3).Value = 'separates'
$ws.Cells.Item(61,4).Value = 'The general approach allows philosophers to ask questions about, for example, what seaprates law from morality, politics, or practical reason.'
$ws.Cells.Item(61,5).Value = 'The general approach allows philosophers to ask questions about, for example, what separates law from morality, politics, or practical reason.'
$ws.Cells.Item(62,2).Value = 'curretn'
$ws.Cells.Item(62,3).Value = 'current'
$ws.Cells.Item(62,4).Value = 'Its curretn constitution was adopted on 7 November 1982 after a constitutional referendum.'
$ws.Cells.Item(62,5).Value = 'Its current constitution was adopted on 7 November 1982 after a constitutional referendum.'
$ws.Cells.Item(63,2).Value = 'assignde'
$ws.Cells.Item(63,3).Value = 'assigned'
$ws.Cells.Item(63,4).Value = 'Justices were required to hold circuit court twice a year in their assignde judicial district.'
$ws.Cells.Item(63,5).Value = 'Justices were required to hold circuit court twice a year in their assigned judicial district.'
$ws.Cells.Item(64,2).Value = 'introdcued'
$ws.Cells.Item(64,3).Value = 'introduced'
$ws.Cells.Item(64,4).Value = 'The idea of sustainable architecture was introdcued in the late 20th century.'
$ws.Cells.Item(64,5).Value = 'The idea of sustainable architecture was introduced in the late 20th century.'
$ws.Cells.Item(65,2).Value = 'Sparatn'
$ws.Cells.Item(65,3).Value = 'Spartan'
$ws.Cells.Item(65,4).Value = 'The Sparatn election of the Ephors, therefore, also predates the reforms of Solon in Athens by approximately 180 years.'
$ws.Cells.Item(65,5).Value = 'The Spartan election of the Ephors, therefore, also predates the reforms of Solon in Athens by approximately 180 years.'
$ws.Cells.Item(66,2).Value = 'eScond'
$ws.Cells.Item(66,3).Value = 'Second'
$ws.Cells.Item(66,4).Value = 'In the 20th century, as a result of scientific progress and the eScond Industrial Revolution, technology stopped being considered a distinct academic discipline and took on the meaning: the systemic use of knowledge to practical ends.'
$ws.Cells.Item(66,5).Value = 'In the 20th century, as a result of scientific progress and the Second Industrial Revolution, technology stopped being considered a distinct academic discipline and took on the meaning: the systemic use of knowledge to practical ends.'
$ws.Cells.Item(67,2).Value = 'amrkets'
$ws.Cells.Item(67,3).Value = 'markets'
$ws.Cells.Item(67,4).Value = 'Trading globally may give consumers and countries the opportunity to be exposed to new amrkets and products.'
$ws.Cells.Item(67,5).Value = 'Trading globally may give consumers and countries the opportunity to be exposed to new markets and products.'
$ws.Cells.Item(68,2).Value = 'expressino'
$ws.Cells.Item(68,3).Value = 'expression'
$ws.Cells.Item(68,4).Value = 'Therefore, freedom of speech and expressino may not be recognized as absolute.'
$ws.Cells.Item(68,5).Value = 'Therefore, freedom of speech and expression may not be recognized as absolute.'
$ws.Cells.Item(69,2).Value = 'hwich'
$ws.Cells.Item(69,3).Value = 'which'
$ws.Cells.Item(69,4).Value = 'Geology is a branch of natural science concerned with the Earth and other astronomical bodies, the rocks of hwich they are composed, and the processes by which they change over time.'
$ws.Cells.Item(69,5).Value = 'Geology is a branch of natural science concerned with the Earth and other astronomical bodies, the rocks of which they are composed, and the processes by which they change over time.'
$ws.Cells.Item(70,2).Value = 'expolres'
$ws.Cells.Item(70,3).Value = 'explores'
$ws.Cells.Item(70,4).Value = 'Archaeology, often referred to as the "anthropology of the past," expolres human activity by examining physical remains.'
$ws.Cells.Item(70,5).Value = 'Archaeology, often referred to as the "anthropology of the past," explores human activity by examining physical remains.'
$ws.Cells.Item(71,2).Value = 'iVtruvian'
$ws.Cells.Item(71,3).Value = 'Vitruvian'
$ws.Cells.Item(71,4).Value = 'The Last Supper is the most reproduced religious painting of all time and his iVtruvian Man drawing is also regarded as a cultural icon.'
$ws.Cells.Item(71,5).Value = 'The Last Supper is the most reproduced religious painting of all time and his Vitruvian Man drawing is also regarded as a cultural icon.'
$ws.Cells.Item(72,2).Value = 'commtitee'
$ws.Cells.Item(72,3).Value = 'committee'
$ws.Cells.Item(72,4).Value = 'To select the commtitee members, a young boy was asked to take out as many leaves as the number of positions available.'
$ws.Cells.Item(72,5).Value = 'To select the committee members, a young boy was asked to take out as many leaves as the number of positions available.'
$ws.Cells.Item(73,2).Value = 'Nteflix'
$ws.Cells.Item(73,3).Value = 'Netflix'
$ws.Cells.Item(73,4).Value = 'Nteflix launched as the first DVD rental and sales website with 30 employees and 925 titles available—nearly all DVDs published.'
$ws.Cells.Item(73,5).Value = 'Netflix launched as the first DVD rental and sales website with 30 employees and 925 titles available—nearly all DVDs published.'
$ws.Cells.Item(74,2).Value = 'rGeece'
$ws.Cells.Item(74,3).Value = 'Greece'
$ws.Cells.Item(74,4).Value = 'Italy, Bulgaria, Romania, rGeece and others entered the war from 1915 onward.'
$ws.Cells.Item(74,5).Value = 'Italy, Bulgaria, Romania, Greece and others entered the war from 1915 onward.'
$ws.Cells.Item(75,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(75,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(75,4).Value = 'It is common for the members of a party to hold similar ideas about politics, and parties may promote specific ideological or policy goals.'
$ws.Cells.Item(75,5).Value = 'It is common for the members of a party to hold similar ideas about politics, and parties may promote specific ideological or policy goals.'
$ws.Cells.Item(76,2).Value = 'bewteen'
$ws.Cells.Item(76,3).Value = 'between'
$ws.Cells.Item(76,4).Value = 'Public memory of 20th-century Marxist-Leninist states has been described as a battleground bewteen anti-communism and anti anti-communism.'
$ws.Cells.Item(76,5).Value = 'Public memory of 20th-century Marxist-Leninist states has been described as a battleground between anti-communism and anti anti-communism.'
$ws.Cells.Item(77,2).Value = 'ebcame'
$ws.Cells.Item(77,3).Value = 'became'
$ws.Cells.Item(77,4).Value = 'He was self-educated and ebcame a lawyer, Illinois state legislator, and U.S.'
$ws.Cells.Item(77,5).Value = 'He was self-educated and became a lawyer, Illinois state legislator, and U.S.'
$ws.Cells.Item(78,2).Value = 'peridos'
$ws.Cells.Item(78,3).Value = 'periods'
$ws.Cells.Item(78,4).Value = 'Six presidents took over the presidency of Egypt after the abolition of the monarchy in 1953, in peridos that included short transitional periods.'
$ws.Cells.Item(78,5).Value = 'Six presidents took over the presidency of Egypt after the abolition of the monarchy in 1953, in periods that included short transitional periods.'
$ws.Cells.Item(79,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(79,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(79,4).Value = 'The linguistic and cultural turns of the mid-20th century, especially, have led to increasingly interpretative, hermeneutic, and philosophical approaches towards the analysis of society.'
$ws.Cells.Item(79,5).Value = 'The linguistic and cultural turns of the mid-20th century, especially, have led to increasingly interpretative, hermeneutic, and philosophical approaches towards the analysis of society.'
$ws.Cells.Item(80,2).Value = 'cpaital'
$ws.Cells.Item(80,3).Value = 'capital'
$ws.Cells.Item(80,4).Value = 'Thenceforth part of the Roman and later Byzantine Empire, the city was the cpaital of the Empire of Trebizond, one of the successor states of the Byzantine Empire after the Fourth Crusade in 1204.'
$ws.Cells.Item(80,5).Value = 'Thenceforth part of the Roman and later Byzantine Empire, the city was the capital of the Empire of Trebizond, one of the successor states of the Byzantine Empire after the Fourth Crusade in 1204.'
$ws.Cells.Item(81,2).Value = 'titel'
$ws.Cells.Item(81,3).Value = 'title'
$ws.Cells.Item(81,4).Value = 'Throughout Argentine history, the office of head of state has undergone many changes, both in its titel as in its features and powers.'
$ws.Cells.Item(81,5).Value = 'Throughout Argentine history, the office of head of state has undergone many changes, both in its title as in its features and powers.'
$ws.Cells.Item(82,2).Value = 'Shakepseare'
$ws.Cells.Item(82,3).Value = 'Shakespeare'
$ws.Cells.Item(82,4).Value = 'Its preface includes a prescient poem by Ben Jonson, a former rival of Shakespeare, who hailed Shakepseare with the now-famous epithet: "not of an age, but for all time".'
$ws.Cells.Item(82,5).Value = 'Its preface includes a prescient poem by Ben Jonson, a former rival of Shakespeare, who hailed Shakespeare with the now-famous epithet: "not of an age, but for all time".'
$ws.Cells.Item(83,2).Value = 'cnetury'
$ws.Cells.Item(83,3).Value = 'century'
$ws.Cells.Item(83,4).Value = 'Having survived the events that caused the fall of the Western Roman Empire in the 5th cnetury AD, it endured until the fall of Constantinople to the Ottoman Empire in 1453.'
$ws.Cells.Item(83,5).Value = 'Having survived the events that caused the fall of the Western Roman Empire in the 5th century AD, it endured until the fall of Constantinople to the Ottoman Empire in 1453.'
$ws.Cells.Item(84,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(84,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(84,4).Value = 'On November 13, 2006, YouTube was purchased by Google for US$1.65 billion (equivalent to $2.39 billion in 2024).'
$ws.Cells.Item(84,5).Value = 'On November 13, 2006, YouTube was purchased by Google for US$1.65 billion (equivalent to $2.39 billion in 2024).'
$ws.Cells.Item(85,2).Value = 'defamaiton'
$ws.Cells.Item(85,3).Value = 'defamation'
$ws.Cells.Item(85,4).Value = 'Excluded were defamaiton of the king''s majesty and the Swedish Church.'
$ws.Cells.Item(85,5).Value = 'Excluded were defamation of the king''s majesty and the Swedish Church.'
$ws.Cells.Item(86,2).Value = 'baolish'
$ws.Cells.Item(86,3).Value = 'abolish'
$ws.Cells.Item(86,4).Value = 'Lincoln is remembered as a martyr and a national hero for his wartime leadership and for his efforts to preserve the Union and baolish slavery.'
$ws.Cells.Item(86,5).Value = 'Lincoln is remembered as a martyr and a national hero for his wartime leadership and for his efforts to preserve the Union and abolish slavery.'
$ws.Cells.Item(87,2).Value = 'Laitn'
$ws.Cells.Item(87,3).Value = 'Latin'
$ws.Cells.Item(87,4).Value = 'In the Spanish-speaking parts of Laitn America the term "federalist" is used in reference to the politics of 19th-century Argentina and Colombia.'
$ws.Cells.Item(87,5).Value = 'In the Spanish-speaking parts of Latin America the term "federalist" is used in reference to the politics of 19th-century Argentina and Colombia.'
$ws.Cells.Item(88,2).Value = 'lPains'
$ws.Cells.Item(88,3).Value = 'Plains'
$ws.Cells.Item(88,4).Value = 'For instance, Winnipeg, a city in the landlocked Great lPains region within Canada, has a January high of −11.3 °C (11.7 °F) and a low of −21.4 °C (−6.5 °F).'
$ws.Cells.Item(88,5).Value = 'For instance, Winnipeg, a city in the landlocked Great Plains region within Canada, has a January high of −11.3 °C (11.7 °F) and a low of −21.4 °C (−6.5 °F).'
$ws.Cells.Item(89,2).Value = 'stduy'
$ws.Cells.Item(89,3).Value = 'study'
$ws.Cells.Item(89,4).Value = 'The term sociology was coined in the late 18th century to describe the scientific stduy of society.'
$ws.Cells.Item(89,5).Value = 'The term sociology was coined in the late 18th century to describe the scientific study of society.'
$ws.Cells.Item(90,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(90,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(90,4).Value = 'It was not until the emergence of the modern evolutionary synthesis from the 1930s to the 1950s that a broad consensus developed in which natural selection was the basic mechanism of evolution.'
$ws.Cells.Item(90,5).Value = 'It was not until the emergence of the modern evolutionary synthesis from the 1930s to the 1950s that a broad consensus developed in which natural selection was the basic mechanism of evolution.'
$ws.Cells.Item(91,2).Value = 'accordnig'
$ws.Cells.Item(91,3).Value = 'according'
$ws.Cells.Item(91,4).Value = 'Sub committees are established accordnig to the issue that the committee receives.'
$ws.Cells.Item(91,5).Value = 'Sub committees are established according to the issue that the committee receives.'
$ws.Cells.Item(92,2).Value = 'supoprted'
$ws.Cells.Item(92,3).Value = 'supported'
$ws.Cells.Item(92,4).Value = 'In domestic policy, Merkel''s Energiewende programme supoprted the development of renewable energy sources and eventually phased out the use of nuclear power in Germany.'
$ws.Cells.Item(92,5).Value = 'In domestic policy, Merkel''s Energiewende programme supported the development of renewable energy sources and eventually phased out the use of nuclear power in Germany.'
$ws.Cells.Item(93,2).Value = 'Prmie'
$ws.Cells.Item(93,3).Value = 'Prime'
$ws.Cells.Item(93,4).Value = 'It was originally a ceremonial post, but became an executive post in 1984 when a new constitution abolished the post of Prmie Minister and transferred its powers to the state president.'
$ws.Cells.Item(93,5).Value = 'It was originally a ceremonial post, but became an executive post in 1984 when a new constitution abolished the post of Prime Minister and transferred its powers to the state president.'
$ws.Cells.Item(94,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(94,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(94,4).Value = 'In 2000, a format change was implemented from a slim quarterly academic journal to a bimonthly magazine.'
$ws.Cells.Item(94,5).Value = 'In 2000, a format change was implemented from a slim quarterly academic journal to a bimonthly magazine.'
$ws.Cells.Item(95,2).Value = 'Caesarae'
$ws.Cells.Item(95,3).Value = 'Caesarea'
$ws.Cells.Item(95,4).Value = 'The only two cities of Cappadocia considered by Strabo to deserve that appellation were Caesarae (originally known as Mazaca) and Tyana, not far from the foot of the Taurus.'
$ws.Cells.Item(95,5).Value = 'The only two cities of Cappadocia considered by Strabo to deserve that appellation were Caesarea (originally known as Mazaca) and Tyana, not far from the foot of the Taurus.'
$ws.Cells.Item(96,2).Value = 'erserve'
$ws.Cells.Item(96,3).Value = 'reserve'
$ws.Cells.Item(96,4).Value = 'The head of state is typically a ceremonial officer, though they may exercise erserve powers to check the Prime Minister in unusual situations.'
$ws.Cells.Item(96,5).Value = 'The head of state is typically a ceremonial officer, though they may exercise reserve powers to check the Prime Minister in unusual situations.'
$ws.Cells.Item(97,2).Value = 'Vcitor'
$ws.Cells.Item(97,3).Value = 'Victor'
$ws.Cells.Item(97,4).Value = 'Presley''s first RCA Vcitor single, "Heartbreak Hotel", was released in January 1956 and became a number-one hit in the US.'
$ws.Cells.Item(97,5).Value = 'Presley''s first RCA Victor single, "Heartbreak Hotel", was released in January 1956 and became a number-one hit in the US.'
$ws.Cells.Item(98,2).Value = 'CRAES'
$ws.Cells.Item(98,3).Value = 'CARES'
$ws.Cells.Item(98,4).Value = 'In response to the COVID-19 pandemic in 2020, he downplayed its severity, contradicted health officials, and signed the CRAES Act.'
$ws.Cells.Item(98,5).Value = 'In response to the COVID-19 pandemic in 2020, he downplayed its severity, contradicted health officials, and signed the CARES Act.'
$ws.Cells.Item(99,2).Value = 'ewbsite'
$ws.Cells.Item(99,3).Value = 'website'
$ws.Cells.Item(99,4).Value = 'The music video for "Infinity Guitars" premiered on NME''s ewbsite on September 19, 2010.'
$ws.Cells.Item(99,5).Value = 'The music video for "Infinity Guitars" premiered on NME''s website on September 19, 2010.'
$ws.Cells.Item(100,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(100,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(100,4).Value = 'The Golden Arches logo was introduced in 1953.'
$ws.Cells.Item(100,5).Value = 'The Golden Arches logo was introduced in 1953.'
$ws.Cells.Item(101,2).Value = 'megadvierse'
$ws.Cells.Item(101,3).Value = 'megadiverse'
$ws.Cells.Item(101,4).Value = 'It is a megadvierse country, with the world''s third-largest land area and third-largest population, exceeding 341 million.'
$ws.Cells.Item(101,5).Value = 'It is a megadiverse country, with the world''s third-largest land area and third-largest population, exceeding 341 million.'
$ws.Cells.Item(102,2).Value = 'Februayr'
$ws.Cells.Item(102,3).Value = 'February'
$ws.Cells.Item(102,4).Value = 'Wilcox, who had purchased 120 acres on Februayr 1, 1887.'
$ws.Cells.Item(102,5).Value = 'Wilcox, who had purchased 120 acres on February 1, 1887.'
$ws.Cells.Item(103,2).Value = 'placse'
$ws.Cells.Item(103,3).Value = 'places'
$ws.Cells.Item(103,4).Value = 'The region is mentioned in the Jewish Mishnah, in Ketubot 13:11, and in several placse in the Talmud, including Yevamot 121a, Hullin 47b.'
$ws.Cells.Item(103,5).Value = 'The region is mentioned in the Jewish Mishnah, in Ketubot 13:11, and in several places in the Talmud, including Yevamot 121a, Hullin 47b.'
$ws.Cells.Item(104,2).Value = 'transmittnig'
$ws.Cells.Item(104,3).Value = 'transmitting'
$ws.Cells.Item(104,4).Value = 'Literature is a method of recording, preserving, and transmittnig knowledge and entertainment.'
$ws.Cells.Item(104,5).Value = 'Literature is a method of recording, preserving, and transmitting knowledge and entertainment.'
$ws.Cells.Item(105,2).Value = 'bsaed'
$ws.Cells.Item(105,3).Value = 'based'
$ws.Cells.Item(105,4).Value = 'Mos Maiorum was a set of rules of conduct bsaed on social norms created over the years by predecessors.'
$ws.Cells.Item(105,5).Value = 'Mos Maiorum was a set of rules of conduct based on social norms created over the years by predecessors.'
$ws.Cells.Item(106,2).Value = 'Ameircan'
$ws.Cells.Item(106,3).Value = 'American'
$ws.Cells.Item(106,4).Value = 'Netflix is an Ameircan subscription video on-demand over-the-top streaming service.'
$ws.Cells.Item(106,5).Value = 'Netflix is an American subscription video on-demand over-the-top streaming service.'
$ws.Cells.Item(107,2).Value = 'dAama'
$ws.Cells.Item(107,3).Value = 'Adama'
$ws.Cells.Item(107,4).Value = 'dAama is a busy transportation center.'
$ws.Cells.Item(107,5).Value = 'Adama is a busy transportation center.'
$ws.Cells.Item(108,2).Value = 'Royla'
$ws.Cells.Item(108,3).Value = 'Royal'
$ws.Cells.Item(108,4).Value = 'He resigned in November 1915 and joined the Royla Scots Fusiliers on the Western Front for six months.'
$ws.Cells.Item(108,5).Value = 'He resigned in November 1915 and joined the Royal Scots Fusiliers on the Western Front for six months.'
$ws.Cells.Item(109,2).Value = 'ibrth'
$ws.Cells.Item(109,3).Value = 'birth'
$ws.Cells.Item(109,4).Value = 'His date of ibrth is unknown but is traditionally observed on 23 April, Saint George''s Day.'
$ws.Cells.Item(109,5).Value = 'His date of birth is unknown but is traditionally observed on 23 April, Saint George''s Day.'
$ws.Cells.Item(110,2).Value = 'Modenr'
$ws.Cells.Item(110,3).Value = 'Modern'
$ws.Cells.Item(110,4).Value = 'Modenr examples include the Yang di-Pertuan Agong (lit.'
$ws.Cells.Item(110,5).Value = 'Modern examples include the Yang di-Pertuan Agong (lit.'
$ws.Cells.Item(111,2).Value = 'presenec'
$ws.Cells.Item(111,3).Value = 'presence'
$ws.Cells.Item(111,4).Value = 'Since then, the ISS has remained continuously inhabited for 25 years and 108 days, the longest continuous human presenec in space.'
$ws.Cells.Item(111,5).Value = 'Since then, the ISS has remained continuously inhabited for 25 years and 108 days, the longest continuous human presence in space.'
$ws.Cells.Item(112,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(112,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(112,4).Value = 'Huntington, have added that while conservatism has core values, going back to the teachings of Edmund Burke, it is also a situational ideology, which aims to conserve diverse social traditions.'
$ws.Cells.Item(112,5).Value = 'Huntington, have added that while conservatism has core values, going back to the teachings of Edmund Burke, it is also a situational ideology, which aims to conserve diverse social traditions.'
$ws.Cells.Item(113,2).Value = 'Aplolo'
$ws.Cells.Item(113,3).Value = 'Apollo'
$ws.Cells.Item(113,4).Value = 'In 1969, Aplolo 11 was the first crewed mission to land on the Moon.'
$ws.Cells.Item(113,5).Value = 'In 1969, Apollo 11 was the first crewed mission to land on the Moon.'
$ws.Cells.Item(114,2).Value = 'sevred'
$ws.Cells.Item(114,3).Value = 'served'
$ws.Cells.Item(114,4).Value = 'A member of the Republican Party, he sevred as the 45th president from 2017 to 2021.'
$ws.Cells.Item(114,5).Value = 'A member of the Republican Party, he served as the 45th president from 2017 to 2021.'
$ws.Cells.Item(115,2).Value = 'menaing'
$ws.Cells.Item(115,3).Value = 'meaning'
$ws.Cells.Item(115,4).Value = 'It was proposed that Kat-patuka came from the Luwian language, menaing "Low Country".'
$ws.Cells.Item(115,5).Value = 'It was proposed that Kat-patuka came from the Luwian language, meaning "Low Country".'
$ws.Cells.Item(116,2).Value = 'etrritorial'
$ws.Cells.Item(116,3).Value = 'territorial'
$ws.Cells.Item(116,4).Value = 'In June 1941, Germany invaded the Soviet Union, opening the Eastern Front and initially making large etrritorial gains along with Axis allies.'
$ws.Cells.Item(116,5).Value = 'In June 1941, Germany invaded the Soviet Union, opening the Eastern Front and initially making large territorial gains along with Axis allies.'
$ws.Cells.Item(117,2).Value = 'amuont'
$ws.Cells.Item(117,3).Value = 'amount'
$ws.Cells.Item(117,4).Value = 'The manifestation of meteorological winter in the northerly snow-prone latitudes is highly variable, depending on elevation, position versus marine winds, and the amuont of precipitation.'
$ws.Cells.Item(117,5).Value = 'The manifestation of meteorological winter in the northerly snow-prone latitudes is highly variable, depending on elevation, position versus marine winds, and the amount of precipitation.'
$ws.Cells.Item(118,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(118,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(118,4).Value = 'In the 19th century, Louis Sullivan declared that "form follows function".'
$ws.Cells.Item(118,5).Value = 'In the 19th century, Louis Sullivan declared that "form follows function".'
$ws.Cells.Item(119,2).Value = 'Shkaespeare'
$ws.Cells.Item(119,3).Value = 'Shakespeare'
$ws.Cells.Item(119,4).Value = 'At the age of 18, Shkaespeare married 26-year-old Anne Hathaway.'
$ws.Cells.Item(119,5).Value = 'At the age of 18, Shakespeare married 26-year-old Anne Hathaway.'
$ws.Cells.Item(120,2).Value = 'peirod'
$ws.Cells.Item(120,3).Value = 'period'
$ws.Cells.Item(120,4).Value = 'The 17th century saw a transformative peirod in British history.'
$ws.Cells.Item(120,5).Value = 'The 17th century saw a transformative period in British history.'
$ws.Cells.Item(121,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(121,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(121,4).Value = 'The UK became the first industrialised country and was the world''s foremost power for the majority of the 19th and early 20th centuries, particularly during the Pax Britannica between 1815 and 1914.'
$ws.Cells.Item(121,5).Value = 'The UK became the first industrialised country and was the world''s foremost power for the majority of the 19th and early 20th centuries, particularly during the Pax Britannica between 1815 and 1914.'
$ws.Cells.Item(122,2).Value = 'histroy'
$ws.Cells.Item(122,3).Value = 'history'
$ws.Cells.Item(122,4).Value = 'The Renaissance (UK: rin-AY-sənss, US: REN-ə-sahnss) is a European period of histroy and cultural movement, very roughly defined as covering the 14th through 17th centuries, though sometimes more narrowly defined for instance as only covering the 15th through 16th centuries.'
$ws.Cells.Item(122,5).Value = 'The Renaissance (UK: rin-AY-sənss, US: REN-ə-sahnss) is a European period of history and cultural movement, very roughly defined as covering the 14th through 17th centuries, though sometimes more narrowly defined for instance as only covering the 15th through 16th centuries.'
$ws.Cells.Item(123,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(123,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(123,4).Value = 'Others explore the physiological and neurobiological processes that underlie cognitive functions and behaviors.'
$ws.Cells.Item(123,5).Value = 'Others explore the physiological and neurobiological processes that underlie cognitive functions and behaviors.'
$ws.Cells.Item(124,2).Value = 'dissovled'
$ws.Cells.Item(124,3).Value = 'dissolved'
$ws.Cells.Item(124,4).Value = 'Because of the Cisplatine War, Rivadavia resigned after a short time, and the office was dissovled shortly thereafter.'
$ws.Cells.Item(124,5).Value = 'Because of the Cisplatine War, Rivadavia resigned after a short time, and the office was dissolved shortly thereafter.'
$ws.Cells.Item(125,2).Value = 'ercourse'
$ws.Cells.Item(125,3).Value = 'recourse'
$ws.Cells.Item(125,4).Value = 'There were no rules of evidence, and little ercourse against an adverse judgment.'
$ws.Cells.Item(125,5).Value = 'There were no rules of evidence, and little recourse against an adverse judgment.'
$ws.Cells.Item(126,2).Value = 'sciencse'
$ws.Cells.Item(126,3).Value = 'sciences'
$ws.Cells.Item(126,4).Value = 'Earth science (also known as geoscience, the geographical sciences or the Earth Sciences), is an all-embracing term for the sciencse related to the planet Earth.'
$ws.Cells.Item(126,5).Value = 'Earth science (also known as geoscience, the geographical sciences or the Earth Sciences), is an all-embracing term for the sciences related to the planet Earth.'
$ws.Cells.Item(127,2).Value = 'membres'
$ws.Cells.Item(127,3).Value = 'members'
$ws.Cells.Item(127,4).Value = 'In 1886, he moved to Paris, where he met membres of the artistic avant-garde, including Émile Bernard and Paul Gauguin, who were seeking new paths beyond Impressionism.'
$ws.Cells.Item(127,5).Value = 'In 1886, he moved to Paris, where he met members of the artistic avant-garde, including Émile Bernard and Paul Gauguin, who were seeking new paths beyond Impressionism.'
$ws.Cells.Item(128,2).Value = 'conglomreate'
$ws.Cells.Item(128,3).Value = 'conglomerate'
$ws.Cells.Item(128,4).Value = 'The Walt Disney Company, commonly known as simply Disney, is an American multinational mass media and entertainment conglomreate headquartered at the Walt Disney Studios complex in Burbank, California.'
$ws.Cells.Item(128,5).Value = 'The Walt Disney Company, commonly known as simply Disney, is an American multinational mass media and entertainment conglomerate headquartered at the Walt Disney Studios complex in Burbank, California.'
$ws.Cells.Item(129,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(129,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(129,4).Value = 'It produces content daily on its website and app, and in four print issues annually.'
$ws.Cells.Item(129,5).Value = 'It produces content daily on its website and app, and in four print issues annually.'
$ws.Cells.Item(130,2).Value = 'stancse'
$ws.Cells.Item(130,3).Value = 'stances'
$ws.Cells.Item(130,4).Value = 'This can lead to problematic stancse which can have large local effects.'
$ws.Cells.Item(130,5).Value = 'This can lead to problematic stances which can have large local effects.'
$ws.Cells.Item(131,2).Value = 'cotninents'
$ws.Cells.Item(131,3).Value = 'continents'
$ws.Cells.Item(131,4).Value = 'The impenetrable defenses enclosed magnificent palaces, domes, and towers, the result of the prosperity Constantinople achieved as the gateway between two cotninents (Europe and Asia) and two seas (the Mediterranean and the Black Sea).'
$ws.Cells.Item(131,5).Value = 'The impenetrable defenses enclosed magnificent palaces, domes, and towers, the result of the prosperity Constantinople achieved as the gateway between two continents (Europe and Asia) and two seas (the Mediterranean and the Black Sea).'
$ws.Cells.Item(132,2).Value = 'strugglde'
$ws.Cells.Item(132,3).Value = 'struggled'
$ws.Cells.Item(132,4).Value = 'It strugglde with internal dissent, especially the Arab Revolt, and engaged in genocide against Armenians, Assyrians, and Greeks.'
$ws.Cells.Item(132,5).Value = 'It struggled with internal dissent, especially the Arab Revolt, and engaged in genocide against Armenians, Assyrians, and Greeks.'
$ws.Cells.Item(133,2).Value = 'evsted'
$ws.Cells.Item(133,3).Value = 'vested'
$ws.Cells.Item(133,4).Value = 'The president is evsted with the "Supreme Executive Power of the Union".'
$ws.Cells.Item(133,5).Value = 'The president is vested with the "Supreme Executive Power of the Union".'
$ws.Cells.Item(134,2).Value = 'histoyr'
$ws.Cells.Item(134,3).Value = 'history'
$ws.Cells.Item(134,4).Value = 'In the histoyr of Europe, the Middle Ages or medieval period lasted approximately from the 5th to the late 15th centuries, comparable with the post-classical period of global history.'
$ws.Cells.Item(134,5).Value = 'In the history of Europe, the Middle Ages or medieval period lasted approximately from the 5th to the late 15th centuries, comparable with the post-classical period of global history.'
$ws.Cells.Item(135,2).Value = 'blmaed'
$ws.Cells.Item(135,3).Value = 'blamed'
$ws.Cells.Item(135,4).Value = 'Austria-Hungary blmaed Serbia, and declared war on 28 July.'
$ws.Cells.Item(135,5).Value = 'Austria-Hungary blamed Serbia, and declared war on 28 July.'
$ws.Cells.Item(136,2).Value = 'salse'
$ws.Cells.Item(136,3).Value = 'sales'
$ws.Cells.Item(136,4).Value = 'When the CD arrived intact, they decided to enter the $16 billion Home-video salse and rental industry.'
$ws.Cells.Item(136,5).Value = 'When the CD arrived intact, they decided to enter the $16 billion Home-video sales and rental industry.'
$ws.Cells.Item(137,2).Value = 'palce'
$ws.Cells.Item(137,3).Value = 'place'
$ws.Cells.Item(137,4).Value = 'All soft landings took palce on the near side of the Moon until January 2019, when Chang''e 4 made the first landing on the far side of the Moon.'
$ws.Cells.Item(137,5).Value = 'All soft landings took place on the near side of the Moon until January 2019, when Chang''e 4 made the first landing on the far side of the Moon.'
$ws.Cells.Item(138,4).Value = 'Foreign Policy was presented as a Gold Winner by the Eddie Awards for "Who Wins in Iraq", in the Consumer News/Commentary/General Interest category.'
$ws.Cells.Item(138,5).Value = 'Foreign Policy was presented as a Gold Winner by the Eddie Awards for "Who Wins in Iraq", in the Consumer News/Commentary/General Interest category.'
$ws.Cells.Item(139,2).Value = 'inlcudes'
$ws.Cells.Item(139,3).Value = 'includes'
$ws.Cells.Item(139,4).Value = 'Criminal law inlcudes the punishment and rehabilitation of people who violate such laws.'
$ws.Cells.Item(139,5).Value = 'Criminal law includes the punishment and rehabilitation of people who violate such laws.'
$ws.Cells.Item(140,2).Value = 'whihc'
$ws.Cells.Item(140,3).Value = 'which'
$ws.Cells.Item(140,4).Value = 'This is partly the effect of states being able to interpret international law in a manner whihc they see fit.'
$ws.Cells.Item(140,5).Value = 'This is partly the effect of states being able to interpret international law in a manner which they see fit.'
$ws.Cells.Item(141,2).Value = 'inclueds'
$ws.Cells.Item(141,3).Value = 'includes'
$ws.Cells.Item(141,4).Value = 'Present-day climate change inclueds both global warming—the ongoing increase in global average temperature—and its wider effects on Earth''s climate system.'
$ws.Cells.Item(141,5).Value = 'Present-day climate change includes both global warming—the ongoing increase in global average temperature—and its wider effects on Earth''s climate system.'
$ws.Cells.Item(142,2).Value = 'Thsee'
$ws.Cells.Item(142,3).Value = 'These'
$ws.Cells.Item(142,4).Value = 'Thsee reforms peaked with the Tanzimat which was the initial reform era of the Ottoman Empire.'
$ws.Cells.Item(142,5).Value = 'These reforms peaked with the Tanzimat which was the initial reform era of the Ottoman Empire.'
$ws.Cells.Item(143,2).Value = 'Thsee'
$ws.Cells.Item(143,3).Value = 'These'
$ws.Cells.Item(143,4).Value = 'Thsee committees are one of auditing tools of the Parliament.'
$ws.Cells.Item(143,5).Value = 'These committees are one of auditing tools of the Parliament.'
$ws.Cells.Item(144,2).Value = 'oriignally'
$ws.Cells.Item(144,3).Value = 'originally'
$ws.Cells.Item(144,4).Value = 'The ISS was oriignally intended to be a laboratory, observatory, and factory while providing transportation, maintenance, and a low Earth orbit staging base for possible future missions to the Moon, Mars, and asteroids.'
$ws.Cells.Item(144,5).Value = 'The ISS was originally intended to be a laboratory, observatory, and factory while providing transportation, maintenance, and a low Earth orbit staging base for possible future missions to the Moon, Mars, and asteroids.'
$ws.Cells.Item(145,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(145,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(145,4).Value = 'In 1968, he returned to the stage in the acclaimed NBC television comeback special Elvis, which led to an extended Las Vegas concert residency and several highly profitable tours.'
$ws.Cells.Item(145,5).Value = 'In 1968, he returned to the stage in the acclaimed NBC television comeback special Elvis, which led to an extended Las Vegas concert residency and several highly profitable tours.'
$ws.Cells.Item(146,2).Value = 'rseponsibility'
$ws.Cells.Item(146,3).Value = 'responsibility'
$ws.Cells.Item(146,4).Value = 'The level of rseponsibility ranged from enacting by-laws about tolls up to the death penalty.'
$ws.Cells.Item(146,5).Value = 'The level of responsibility ranged from enacting by-laws about tolls up to the death penalty.'
$ws.Cells.Item(147,2).Value = 'closley'
$ws.Cells.Item(147,3).Value = 'closely'
$ws.Cells.Item(147,4).Value = 'Legal history is closley connected to the development of civilizations and operates in the wider context of social history.'
$ws.Cells.Item(147,5).Value = 'Legal history is closely connected to the development of civilizations and operates in the wider context of social history.'
$ws.Cells.Item(148,2).Value = 'ubildings'
$ws.Cells.Item(148,3).Value = 'buildings'
$ws.Cells.Item(148,4).Value = 'Its ubildings and facilities are scattered throughout the city centre and around the town.'
$ws.Cells.Item(148,5).Value = 'Its buildings and facilities are scattered throughout the city centre and around the town.'
$ws.Cells.Item(149,2).Value = 'Ltierature'
$ws.Cells.Item(149,3).Value = 'Literature'
$ws.Cells.Item(149,4).Value = 'He was awarded the Nobel Prize in Ltierature in 1953.'
$ws.Cells.Item(149,5).Value = 'He was awarded the Nobel Prize in Literature in 1953.'
$ws.Cells.Item(150,2).Value = 'Intrenational'
$ws.Cells.Item(150,3).Value = 'International'
$ws.Cells.Item(150,4).Value = 'Intrenational law, also known as public international law and the law of nations, is the set of rules, norms, legal customs and standards that states and other actors feel an obligation to, and generally do, obey in their mutual relations.'
$ws.Cells.Item(150,5).Value = 'International law, also known as public international law and the law of nations, is the set of rules, norms, legal customs and standards that states and other actors feel an obligation to, and generally do, obey in their mutual relations.'
$ws.Cells.Item(151,2).Value = 'rehabiltiative'
$ws.Cells.Item(151,3).Value = 'rehabilitative'
$ws.Cells.Item(151,4).Value = 'Criminal procedure is a formalized official activity that authenticates the fact of commission of a crime and authorizes punitive or rehabiltiative treatment of the offender.'
$ws.Cells.Item(151,5).Value = 'Criminal procedure is a formalized official activity that authenticates the fact of commission of a crime and authorizes punitive or rehabilitative treatment of the offender.'
$ws.Cells.Item(152,2).Value = 'davances'
$ws.Cells.Item(152,3).Value = 'advances'
$ws.Cells.Item(152,4).Value = 'Other technological davances made during the Paleolithic era include clothing and shelter.'
$ws.Cells.Item(152,5).Value = 'Other technological advances made during the Paleolithic era include clothing and shelter.'
$ws.Cells.Item(153,2).Value = 'Serivces'
$ws.Cells.Item(153,3).Value = 'Services'
$ws.Cells.Item(153,4).Value = 'Serivces are also traded, such as in tourism, banking, consulting, and transportation.'
$ws.Cells.Item(153,5).Value = 'Services are also traded, such as in tourism, banking, consulting, and transportation.'
$ws.Cells.Item(154,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(154,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(154,4).Value = 'A Moon landing or lunar landing is the arrival of a spacecraft on the surface of the Moon, including both crewed and robotic missions.'
$ws.Cells.Item(154,5).Value = 'A Moon landing or lunar landing is the arrival of a spacecraft on the surface of the Moon, including both crewed and robotic missions.'
$ws.Cells.Item(155,2).Value = 'trnasacted'
$ws.Cells.Item(155,3).Value = 'transacted'
$ws.Cells.Item(155,4).Value = 'Economic transactions occur when two groups or parties agree to the value or price of the trnasacted good or service, commonly expressed in a certain currency.'
$ws.Cells.Item(155,5).Value = 'Economic transactions occur when two groups or parties agree to the value or price of the transacted good or service, commonly expressed in a certain currency.'
$ws.Cells.Item(156,2).Value = 'histroy'
$ws.Cells.Item(156,3).Value = 'history'
$ws.Cells.Item(156,4).Value = 'In a more general sense, the term histroy refers not to an academic field but to the past itself, times in the past, or to individual texts about the past.'
$ws.Cells.Item(156,5).Value = 'In a more general sense, the term history refers not to an academic field but to the past itself, times in the past, or to individual texts about the past.'
$ws.Cells.Item(157,2).Value = 'Vicni'
$ws.Cells.Item(157,3).Value = 'Vinci'
$ws.Cells.Item(157,4).Value = 'It saw myriad artistic developments and contributions from such polymaths as Leonardo da Vicni and Michelangelo, who inspired the term "Renaissance man".'
$ws.Cells.Item(157,5).Value = 'It saw myriad artistic developments and contributions from such polymaths as Leonardo da Vinci and Michelangelo, who inspired the term "Renaissance man".'
$ws.Cells.Item(158,2).Value = 'psychologisst'
$ws.Cells.Item(158,3).Value = 'psychologists'
$ws.Cells.Item(158,4).Value = 'Typically the latter group of psychologisst work in academic settings (e.g., universities, medical schools, or hospitals).'
$ws.Cells.Item(158,5).Value = 'Typically the latter group of psychologists work in academic settings (e.g., universities, medical schools, or hospitals).'
$ws.Cells.Item(159,2).Value = 'Parliamentayr'
$ws.Cells.Item(159,3).Value = 'Parliamentary'
$ws.Cells.Item(159,4).Value = 'Most of Ireland seceded from the UK in 1922 as the Irish Free State, and the Royal and Parliamentayr Titles Act 1927 created the present United Kingdom.'
$ws.Cells.Item(159,5).Value = 'Most of Ireland seceded from the UK in 1922 as the Irish Free State, and the Royal and Parliamentary Titles Act 1927 created the present United Kingdom.'
$ws.Cells.Item(160,2).Value = 'oscialist'
$ws.Cells.Item(160,3).Value = 'socialist'
$ws.Cells.Item(160,4).Value = 'Communism is a part of the broader oscialist movement.'
$ws.Cells.Item(160,5).Value = 'Communism is a part of the broader socialist movement.'
$ws.Cells.Item(161,2).Value = 'beacme'
$ws.Cells.Item(161,3).Value = 'became'
$ws.Cells.Item(161,4).Value = 'During the reigns of Selim I and Suleiman the Magnificent, the Ottoman Empire beacme a global power.'
$ws.Cells.Item(161,5).Value = 'During the reigns of Selim I and Suleiman the Magnificent, the Ottoman Empire became a global power.'
$ws.Cells.Item(162,2).Value = 'nuiversity'
$ws.Cells.Item(162,3).Value = 'university'
$ws.Cells.Item(162,4).Value = 'In 1231, 22 years after its founding, the nuiversity was recognised with a royal charter, granted by King Henry III.'
$ws.Cells.Item(162,5).Value = 'In 1231, 22 years after its founding, the university was recognised with a royal charter, granted by King Henry III.'
$ws.Cells.Item(163,2).Value = 'hitsory'
$ws.Cells.Item(163,3).Value = 'history'
$ws.Cells.Item(163,4).Value = 'The study of books and other texts as artifacts or traditions is instead encompassed by textual criticism or the hitsory of the book.'
$ws.Cells.Item(163,5).Value = 'The study of books and other texts as artifacts or traditions is instead encompassed by textual criticism or the history of the book.'
$ws.Cells.Item(164,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(164,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(164,4).Value = 'A product that is transferred or sold from a party in one country to a party in another country is an export from the originating country, and an import to the country receiving that product.'
$ws.Cells.Item(164,5).Value = 'A product that is transferred or sold from a party in one country to a party in another country is an export from the originating country, and an import to the country receiving that product.'
$ws.Cells.Item(165,1).Value = 'transposition'
$ws.Cells.Item(165,2).Value = 'McoDnald'
$ws.Cells.Item(165,3).Value = 'McDonald'
$ws.Cells.Item(165,4).Value = 'Clown mascot Ronald McoDnald was introduced in 1963 to market the chain to children.'
$ws.Cells.Item(165,5).Value = 'Clown mascot Ronald McDonald was introduced in 1963 to market the chain to children.'
$ws.Cells.Item(166,1).Value = 'transposition'
$ws.Cells.Item(166,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(166,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(166,4).Value = 'Cilicia was the name given to the district in which Caesarea, the capital of the whole country, was situated.'
$ws.Cells.Item(166,5).Value = 'Cilicia was the name given to the district in which Caesarea, the capital of the whole country, was situated.'
$ws.Cells.Item(167,1).Value = 'transposition'
$ws.Cells.Item(167,2).Value = 'oFlio'
$ws.Cells.Item(167,3).Value = 'Folio'
$ws.Cells.Item(167,4).Value = 'oFlio Magazine Gold Editorial Excellence (Eddie) Award – Consumer Magazine, News/Commentary/General Interest (single article), "What America Must Do" by Kenneth Rogoff, Jan/Feb 2008.'
$ws.Cells.Item(167,5).Value = 'Folio Magazine Gold Editorial Excellence (Eddie) Award – Consumer Magazine, News/Commentary/General Interest (single article), "What America Must Do" by Kenneth Rogoff, Jan/Feb 2008.'
$ws.Cells.Item(168,1).Value = 'transposition'
$ws.Cells.Item(168,2).Value = 'knonw'
$ws.Cells.Item(168,3).Value = 'known'
$ws.Cells.Item(168,4).Value = 'The Western Front is knonw as the Greco-Turkish War.'
$ws.Cells.Item(168,5).Value = 'The Western Front is known as the Greco-Turkish War.'
$ws.Cells.Item(169,1).Value = 'transposition'
$ws.Cells.Item(169,2).Value = 'Knigdom'
$ws.Cells.Item(169,3).Value = 'Kingdom'
$ws.Cells.Item(169,4).Value = 'The Acts of Union 1707 declared that the Knigdom of England and the Kingdom of Scotland were "United into One Kingdom by the Name of Great Britain".'
$ws.Cells.Item(169,5).Value = 'The Acts of Union 1707 declared that the Kingdom of England and the Kingdom of Scotland were "United into One Kingdom by the Name of Great Britain".'
$ws.Cells.Item(170,1).Value = 'transposition'
$ws.Cells.Item(170,2).Value = 'paplications'
$ws.Cells.Item(170,3).Value = 'applications'
$ws.Cells.Item(170,4).Value = 'The paplications of various fields of chemistry are used frequently for economic purposes in the chemical industry.'
$ws.Cells.Item(170,5).Value = 'The applications of various fields of chemistry are used frequently for economic purposes in the chemical industry.'
$ws.Cells.Item(171,1).Value = 'transposition'
$ws.Cells.Item(171,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(171,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(171,4).Value = 'With moisture being trapped, the local climate thus has high winter rainfall, while the interior bay setting results in very hot summers for a coastal city.'
$ws.Cells.Item(171,5).Value = 'With moisture being trapped, the local climate thus has high winter rainfall, while the interior bay setting results in very hot summers for a coastal city.'
$ws.Cells.Item(172,1).Value = 'transposition'
$ws.Cells.Item(172,2).Value = 'proir'
$ws.Cells.Item(172,3).Value = 'prior'
$ws.Cells.Item(172,4).Value = 'Be a resident of Mexico for the entire year proir to the election (although absences of 30 days or fewer are explicitly stated not to interrupt residency).'
$ws.Cells.Item(172,5).Value = 'Be a resident of Mexico for the entire year prior to the election (although absences of 30 days or fewer are explicitly stated not to interrupt residency).'
$ws.Cells.Item(173,1).Value = 'transposition'
$ws.Cells.Item(173,2).Value = 'Auugst'
$ws.Cells.Item(173,3).Value = 'August'
$ws.Cells.Item(173,4).Value = 'As of Auugst 2025, 290 individuals from 26 countries had visited the station.'
$ws.Cells.Item(173,5).Value = 'As of August 2025, 290 individuals from 26 countries had visited the station.'
$ws.Cells.Item(174,1).Value = 'transposition'
$ws.Cells.Item(174,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(174,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(174,4).Value = 'Ancient Athens was the cradle of democracy.'
$ws.Cells.Item(174,5).Value = 'Ancient Athens was the cradle of democracy.'
$ws.Cells.Item(175,1).Value = 'transposition'
$ws.Cells.Item(175,2).Value = 'yaers'
$ws.Cells.Item(175,3).Value = 'years'
$ws.Cells.Item(175,4).Value = 'Suffrage is universal for citizens 18 yaers of age and older.'
$ws.Cells.Item(175,5).Value = 'Suffrage is universal for citizens 18 years of age and older.'
$ws.Cells.Item(176,1).Value = 'transposition'
$ws.Cells.Item(176,2).Value = 'toawrd'
$ws.Cells.Item(176,3).Value = 'toward'
$ws.Cells.Item(176,4).Value = 'Conversely, winter in the Southern Hemisphere occurs when the Northern Hemisphere is tilted more toawrd the Sun.'
$ws.Cells.Item(176,5).Value = 'Conversely, winter in the Southern Hemisphere occurs when the Northern Hemisphere is tilted more toward the Sun.'
$ws.Cells.Item(177,1).Value = 'transposition'
$ws.Cells.Item(177,2).Value = 'Tespit Edilemedi'
$ws.Cells.Item(177,3).Value = 'Tespit Edilemedi'
$ws.Cells.Item(177,4).Value = 'Fearing more violence from Oxford townsfolk, University of Oxford scholars began leaving Oxford for more hospitable cities, including Paris, Reading, and Cambridge.'
$ws.Cells.Item(177,5).Value = 'Fearing more violence from Oxford townsfolk, University of Oxford scholars began leaving Oxford for more hospitable cities, including Paris, Reading, and Cambridge.'
$ws.Cells.Item(178,1).Value = 'transposition'
$ws.Cells.Item(178,2).Value = 'distnict'
$ws.Cells.Item(178,3).Value = 'distinct'
$ws.Cells.Item(178,4).Value = 'Some scientists include as part of the spheres of the Earth, the cryosphere (corresponding to ice) as a distnict portion of the hydrosphere, as well as the pedosphere (to soil) as an active and intermixed sphere.'
$ws.Cells.Item(178,5).Value = 'Some scientists include as part of the spheres of the Earth, the cryosphere (corresponding to ice) as a distinct portion of the hydrosphere, as well as the pedosphere (to soil) as an active and intermixed sphere.'
$ws.Cells.Item(179,1).Value = 'transposition'
$ws.Cells.Item(179,2).Value = 'Ottoamn'
$ws.Cells.Item(179,3).Value = 'Ottoman'
$ws.Cells.Item(179,4).Value = 'When the millet system started to lose its efficiency due to the rise of nationalism within its borders, the Ottoamn Empire explored new ways of governing its territory composed of diverse populations.'
$ws.Cells.Item(179,5).Value = 'When the millet system started to lose its efficiency due to the rise of nationalism within its borders, the Ottoman Empire explored new ways of governing its territory composed of diverse populations.'
$ws.Cells.Item(180,1).Value = 'transposition'
$ws.Cells.Item(180,2).Value = 'Anegla'
$ws.Cells.Item(180,3).Value = 'Angela'
$ws.Cells.Item(180,4).Value = 'Anegla Dorothea Merkel (née Kasner; born 17 July 1954) is a German retired politician who served as Chancellor of Germany from 2005 to 2021.'
$ws.Cells.Item(180,5).Value = 'Angela Dorothea Merkel (née Kasner; born 17 July 1954) is a German retired politician who served as Chancellor of Germany from 2005 to 2021.'
